# Slide 7 ("Warning: Turn Off Tabs!!") - remove the NotePad++ / TextWrangler
# bullet points and split the leading "Python " word off of the following
# bullet's first run (accessibility/arrow-color cleanup commit).

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(7)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

# Paragraph 4 is "NotePad++:  Settings -> Preferences -> Language Menu/Tab Settings"
$tr.Paragraphs(4, 1).Delete()

# After the delete above, paragraph 4 is now
# "TextWrangler:  TextWrangler -> Preferences -> Editor Defaults"
$tr.Paragraphs(4, 1).Delete()

# Paragraph 4 is now "Python cares a *lot* about how far a line is indented. ..."
# Split its first run so "Python " becomes its own run, leaving
# "cares a *lot* about how far a line is " as the following run.
$target = $tr.Paragraphs(4, 1)
$lead = $tr.Characters($target.Start, 7)
$lead.Text = "Python "
